$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Taille" (height) column (E) for the players that were missing it.
$ws.Range("E2").Value  = "1m80"
$ws.Range("E3").Value  = "1m87"
$ws.Range("E4").Value  = "1m84"
$ws.Range("E5").Value  = "1m81"
$ws.Range("E7").Value  = "1m91"
$ws.Range("E11").Value = "1m83"
$ws.Range("E12").Value = "1m89"
$ws.Range("E13").Value = "1m72"
$ws.Range("E16").Value = "1m93"
$ws.Range("E17").Value = "1m74"
$ws.Range("E22").Value = "1m69"
$ws.Range("E25").Value = "1m85"

# Realign the first player's name cell so it matches the rest of the column
# (center / center instead of left / top).
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108

# Put the selection where the author last left it.
[void]$ws.Range("F26").Select()
